# Weekly price-sheet update: insert a new daily record as row 187,
# pushing the existing rows 187-258 down to 188-259 (dimension grows to R259).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187; rows below shift down automatically,
# and the new row inherits formatting (date style, etc.) from the row above.
$ws.Rows.Item(187).Insert()

$ws.Range('A187').Value = 9
$ws.Range('B187').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C187').Value = 'Metropolitana'
$ws.Range('D187').Value = 44627
$ws.Range('E187').Value = 13
$ws.Range('F187').Value = 100112021
$ws.Range('G187').Value = 'Ají'
$ws.Range('H187').Value = 'Americana (o)'
$ws.Range('I187').Value = 'Primera'
$ws.Range('J187').Value = 61
$ws.Range('K187').Value = 1000
$ws.Range('L187').Value = 11000
$ws.Range('M187').Value = 6082
$ws.Range('N187').Value = '$/caja 25 kilos'
$ws.Range('O187').Value = 'Provincia de Limarí'
$ws.Range('P187').Value = 243
$ws.Range('Q187').Value = 25
$ws.Range('R187').Value = 'Hortaliza'
